$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '315.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.59%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '39.41'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-0.75%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.139'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.52%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08194'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '1.06%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.990'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2.25%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.367'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '3.08%'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '8.320'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.06%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9380'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.31%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1304'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-8.18%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1973'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.77%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09091'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.14%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03479'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-1.09%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09760'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.41%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001409'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.74%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006373'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '8.86%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.635'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-7.62%'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.27%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1316'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.37%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.957'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '6.32%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2489'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2.77%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04364'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.08%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001243'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.03%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004772'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '9.42%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003889'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '198.79%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-7.56%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02242'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '9.83%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05209'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2.98%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007744'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '4.67%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01034'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '5.54%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1400'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2.64%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002100'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-1.56%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.008885'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-6.30%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006818'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '7.16%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.04%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002992'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '9.70%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.001690'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '29.93%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002100'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.04%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002000'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.04%'
